$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "71.056.08"
$ws.Range("E2").Value = "  +6.55%  "
$ws.Range("D3").Value = "3.656.26"
$ws.Range("E3").Value = "  +18.47%  "
$ws.Range("E4").Value = "  +0.08%  "
$c = $ws.Range("D5")
$c.NumberFormat = "@"
$c.Value = "616.97"
$c.Style = "Normal"
$ws.Range("E5").Value = "  +7.17%  "
$c = $ws.Range("D6")
$c.NumberFormat = "@"
$c.Value = "180.66"
$c.Style = "Normal"
$ws.Range("E6").Value = "  +1.96%  "
$ws.Range("D7").Value = "3.654.10"
$ws.Range("E7").Value = "  +18.39%  "
$ws.Range("E8").Value = "  +0.00%  "
$c = $ws.Range("D9")
$c.NumberFormat = "@"
$c.Value = "0.535"
$c.Style = "Normal"
$ws.Range("E9").Value = "  +4.38%  "
$c = $ws.Range("D10")
$c.NumberFormat = "@"
$c.Value = "0.162"
$c.Style = "Normal"
$ws.Range("E10").Value = "  +8.04%  "
$c = $ws.Range("D11")
$c.NumberFormat = "@"
$c.Value = "6.62"
$c.Style = "Normal"
$ws.Range("E11").Value = "  +4.34%  "
$c = $ws.Range("D12")
$c.NumberFormat = "@"
$c.Value = "0.498"
$c.Style = "Normal"
$ws.Range("E12").Value = "  +7.14%  "
$c = $ws.Range("D13")
$c.NumberFormat = "@"
$c.Value = "40.50"
$c.Style = "Normal"
$ws.Range("E13").Value = "  +12.99%  "
$c = $ws.Range("D14")
$c.NumberFormat = "@"
$c.Value = "0.0000254"
$c.Style = "Normal"
$ws.Range("E14").Value = "  +6.16%  "
$ws.Range("D15").Value = "4.264.45"
$ws.Range("E15").Value = "  +18.56%  "
$ws.Range("D16").Value = "71.052.12"
$ws.Range("E16").Value = "  +6.68%  "
$ws.Range("D17").Value = "3.655.58"
$ws.Range("E17").Value = "  +18.58%  "
$ws.Range("E18").Value = "  +1.27%  "
$c = $ws.Range("D19")
$c.NumberFormat = "@"
$c.Value = "7.51"
$c.Style = "Normal"
$ws.Range("E19").Value = "  +8.05%  "
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$c = $ws.Range("D20")
$c.NumberFormat = "@"
$c.Value = "519.90"
$c.Style = "Normal"
$ws.Range("E20").Value = "  +8.37%  "
$ws.Range("B21").Value = "Chainlink"
$ws.Range("C21").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$c = $ws.Range("D21")
$c.NumberFormat = "@"
$c.Value = "16.87"
$c.Style = "Normal"
$ws.Range("E21").Value = "  +0.43%  "
$c = $ws.Range("D22")
$c.NumberFormat = "@"
$c.Value = "9.27"
$c.Style = "Normal"
$ws.Range("E22").Value = "  +20.18%  "
$c = $ws.Range("D23")
$c.NumberFormat = "@"
$c.Value = "0.742"
$c.Style = "Normal"
$ws.Range("E23").Value = "  +8.19%  "
$c = $ws.Range("D24")
$c.NumberFormat = "@"
$c.Value = "88.30"
$c.Style = "Normal"
$ws.Range("E24").Value = "  +6.13%  "
$c = $ws.Range("D25")
$c.NumberFormat = "@"
$c.Value = "2.46"
$c.Style = "Normal"
$ws.Range("E25").Value = "  +11.45%  "
$c = $ws.Range("D26")
$c.NumberFormat = "@"
$c.Value = "13.42"
$c.Style = "Normal"
$ws.Range("E26").Value = "  +6.51%  "
$c = $ws.Range("D27")
$c.NumberFormat = "@"
$c.Value = "10.92"
$c.Style = "Normal"
$ws.Range("E27").Value = "  +8.61%  "
$c = $ws.Range("D28")
$c.NumberFormat = "@"
$c.Value = "1.00"
$c.Style = "Normal"
$ws.Range("E28").Value = "  -0.15%  "
$c = $ws.Range("D29")
$c.NumberFormat = "@"
$c.Value = "2.54"
$c.Style = "Normal"
$ws.Range("E29").Value = "  +11.19%  "
$c = $ws.Range("D30")
$c.NumberFormat = "@"
$c.Value = "8.09"
$c.Style = "Normal"
$ws.Range("E30").Value = "  +1.92%  "
$c = $ws.Range("D31")
$c.NumberFormat = "@"
$c.Value = "2.88"
$c.Style = "Normal"
$ws.Range("E31").Value = "  +11.18%  "
$ws.Range("B32").Value = "EthereumClassic"
$ws.Range("C32").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$c = $ws.Range("D32")
$c.NumberFormat = "@"
$c.Value = "31.65"
$c.Style = "Normal"
$ws.Range("E32").Value = "  +13.53%  "
$ws.Range("B33").Value = "PEPE"
$ws.Range("C33").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$c = $ws.Range("D33")
$c.NumberFormat = "@"
$c.Value = "0.0000110"
$c.Style = "Normal"
$ws.Range("E33").Value = "  +17.28%  "
$ws.Range("E34").Value = "  +3.96%  "
$ws.Range("E35").Value = "  +0.11%  "
$c = $ws.Range("D36")
$c.NumberFormat = "@"
$c.Value = "6.13"
$c.Style = "Normal"
$ws.Range("E36").Value = "  +10.34%  "
$c = $ws.Range("D37")
$c.NumberFormat = "@"
$c.Value = "1.01"
$c.Style = "Normal"
$ws.Range("E37").Value = "  +8.19%  "
$c = $ws.Range("D38")
$c.NumberFormat = "@"
$c.Value = "0.345"
$c.Style = "Normal"
$ws.Range("E38").Value = "  +12.12%  "
$c = $ws.Range("D39")
$c.NumberFormat = "@"
$c.Value = "2.18"
$c.Style = "Normal"
$ws.Range("E39").Value = "  +10.41%  "
$c = $ws.Range("D40")
$c.NumberFormat = "@"
$c.Value = "51.62"
$c.Style = "Normal"
$ws.Range("E40").Value = "  +5.50%  "
$ws.Range("E41").Value = "  +5.46%  "
$c = $ws.Range("D42")
$c.NumberFormat = "@"
$c.Value = "45.49"
$c.Style = "Normal"
$ws.Range("E42").Value = "  -5.07%  "
$ws.Range("B43").Value = "Maker"
$ws.Range("C43").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D43").Value = "3.120.95"
$ws.Range("E43").Value = "  +12.22%  "
$ws.Range("B44").Value = "Cosmos"
$ws.Range("C44").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$c = $ws.Range("D44")
$c.NumberFormat = "@"
$c.Value = "8.79"
$c.Style = "Normal"
$ws.Range("E44").Value = "  +5.94%  "
$c = $ws.Range("D45")
$c.NumberFormat = "@"
$c.Value = "419.69"
$c.Style = "Normal"
$ws.Range("E45").Value = "  +13.13%  "
$c = $ws.Range("D46")
$c.NumberFormat = "@"
$c.Value = "2.76"
$c.Style = "Normal"
$ws.Range("E46").Value = "  +3.93%  "
$c = $ws.Range("D47")
$c.NumberFormat = "@"
$c.Value = "28.84"
$c.Style = "Normal"
$ws.Range("E47").Value = "  +16.88%  "
$ws.Range("E48").Value = "  +7.91%  "
$c = $ws.Range("D49")
$c.NumberFormat = "@"
$c.Value = "138.99"
$c.Style = "Normal"
$ws.Range("E49").Value = "  +2.64%  "
$c = $ws.Range("D51")
$c.NumberFormat = "@"
$c.Value = "2.46"
$c.Style = "Normal"
$ws.Range("E51").Value = "  +11.62%  "
